# Update seed data for purchase receive header sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (P REV 01 / test-id-14 / P REQ 01 / RECEIVED (NO ISSUE)):
#   only the RECEIVER date moves from 2025-01-08 (45665) to 2025-01-05 (45662).
$ws.Range("B2").Value = 45662

# Row 3 (P REV 02 / test-id-14 / ... / RECEIVED (NO ISSUE)):
#   date moves the same way, and the HEADER CODE changes from P REQ 01 to P REQ 02.
$ws.Range("B3").Value = 45662
$ws.Range("D3").Value = "P REQ 02"

# Row 4 (P REV 03 / test-id-15 / P REQ 02 / RECEIVED (WITH ISSUE)) is removed entirely.
$ws.Rows("4:4").Delete()

# Selection ends up on E4 (just past the now-smaller used range).
$ws.Range("E4").Select()
